$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.4061527774565
$ws.Range("C2").Value = 14.11593552726775
$ws.Range("E2").Value = 17.37733260090663
$ws.Range("F2").Value = 35.75600207900801
$ws.Range("G2").Value = 27.20488201108558
$ws.Range("H2").Value = 14.00434511122224
$ws.Range("J2").Value = 7.648470302336461
$ws.Range("L2").Value = 12.68603526850549
$ws.Range("M2").Value = 16.38472673520936
$ws.Range("O2").Value = 21.09027059365978
$ws.Range("B3").Value = 14.83192185359091
$ws.Range("C3").Value = 14.04946553717104
$ws.Range("E3").Value = 17.4328929074303
$ws.Range("F3").Value = 35.87152619303483
$ws.Range("G3").Value = 27.3575832323414
$ws.Range("H3").Value = 14.06670789611541
$ws.Range("J3").Value = 7.633912745987269
$ws.Range("L3").Value = 12.67226545424522
$ws.Range("M3").Value = 16.23766372181362
$ws.Range("O3").Value = 21.2009953745339
$ws.Range("B4").Value = 14.46813335856349
$ws.Range("C4").Value = 14.00953775688094
$ws.Range("E4").Value = 17.4692702855032
$ws.Range("F4").Value = 35.95160492497351
$ws.Range("G4").Value = 27.46202823575716
$ws.Range("H4").Value = 14.10755653212881
$ws.Range("J4").Value = 7.625059102307712
$ws.Range("L4").Value = 12.66511736767296
$ws.Range("M4").Value = 16.14804614785492
$ws.Range("O4").Value = 21.27426615550019
$ws.Range("B5").Value = 14.31726790073916
$ws.Range("C5").Value = 13.9934991597292
$ws.Range("E5").Value = 17.484664234661
$ws.Range("F5").Value = 35.9865314651848
$ws.Range("G5").Value = 27.50726054477727
$ws.Range("H5").Value = 14.12484594775397
$ws.Range("J5").Value = 7.621473400852542
$ws.Range("L5").Value = 12.66253527196791
$ws.Range("M5").Value = 16.11172606857212
$ws.Range("O5").Value = 21.30545099950261
$ws.Range("B6").Value = 14.29206491240742
$ws.Range("C6").Value = 13.99085027295684
$ws.Range("E6").Value = 17.48725483046383
$ws.Range("F6").Value = 35.99246934401761
$ws.Range("G6").Value = 27.51493210279419
$ws.Range("H6").Value = 14.12775569780456
$ws.Range("J6").Value = 7.620879367039016
$ws.Range("L6").Value = 12.6621265598637
$ws.Range("M6").Value = 16.10570805190038
$ws.Range("O6").Value = 21.31070925446357
$ws.Range("B7").Value = 14.46610905718371
$ws.Range("C7").Value = 14.00932050150332
$ws.Range("E7").Value = 17.46947558523466
$ws.Range("F7").Value = 35.95206667737312
$ws.Range("G7").Value = 27.46262746590378
$ws.Range("H7").Value = 14.10778709850168
$ws.Range("J7").Value = 7.625010653077026
$ws.Range("L7").Value = 12.66508120233235
$ws.Range("M7").Value = 16.14755547566285
$ws.Range("O7").Value = 21.27468135810985
$ws.Range("B8").Value = 15.21059738061076
$ws.Range("C8").Value = 14.09283889000743
$ws.Range("E8").Value = 17.39602055910726
$ws.Range("F8").Value = 35.79393286199828
$ws.Range("G8").Value = 27.25530548307415
$ws.Range("H8").Value = 14.02531703704963
$ws.Range("J8").Value = 7.643433818347319
$ws.Range("L8").Value = 12.68101762616362
$ws.Range("M8").Value = 16.33389576665963
$ws.Range("O8").Value = 21.12734994532261
$ws.Range("B9").Value = 16.57390991951987
$ws.Range("C9").Value = 14.26320504067749
$ws.Range("E9").Value = 17.26989655775607
$ws.Range("F9").Value = 35.55666564688663
$ws.Range("G9").Value = 26.93427361113246
$ws.Range("H9").Value = 13.88388132974649
$ws.Range("J9").Value = 7.680194365159132
$ws.Range("L9").Value = 12.72253401852412
$ws.Range("M9").Value = 16.70342102554571
$ws.Range("O9").Value = 20.88048806264947
$ws.Range("B10").Value = 17.50773938986412
$ws.Range("C10").Value = 14.39179479597714
$ws.Range("E10").Value = 17.18810874588744
$ws.Range("F10").Value = 35.42707821378682
$ws.Range("G10").Value = 26.75153155108829
$ws.Range("H10").Value = 13.79232831837282
$ws.Range("J10").Value = 7.707535715774289
$ws.Range("L10").Value = 12.75915093247197
$ws.Range("M10").Value = 16.97569512065611
$ws.Range("O10").Value = 20.72491824639261
$ws.Range("B11").Value = 17.91633022079394
$ws.Range("C11").Value = 14.45090675629091
$ws.Range("E11").Value = 17.1532524717835
$ws.Range("F11").Value = 35.37790096485559
$ws.Range("G11").Value = 26.68013179375797
$ws.Range("H11").Value = 13.7533610474931
$ws.Range("J11").Value = 7.720035015157698
$ws.Range("L11").Value = 12.77710355708028
$ws.Range("M11").Value = 17.09936376950794
$ws.Range("O11").Value = 20.65978314902862
$ws.Range("B12").Value = 18.06861716882271
$ws.Range("C12").Value = 14.47336824154272
$ws.Range("E12").Value = 17.14039039634902
$ws.Range("F12").Value = 35.36068844264499
$ws.Range("G12").Value = 26.65479698103813
$ws.Range("H12").Value = 13.73899065867806
$ws.Range("J12").Value = 7.724776002003317
$ws.Range("L12").Value = 12.78408495942314
$ws.Range("M12").Value = 17.14613596888228
$ws.Range("O12").Value = 20.63593149197057
$ws.Range("B13").Value = 18.03592939726093
$ws.Range("C13").Value = 14.4685275150195
$ws.Range("E13").Value = 17.14314548511341
$ws.Range("F13").Value = 35.36433269801251
$ws.Range("G13").Value = 26.66017731129526
$ws.Range("H13").Value = 13.7420684228089
$ws.Range("J13").Value = 7.723754617313077
$ws.Range("L13").Value = 12.78257329610082
$ws.Range("M13").Value = 17.13606582857831
$ws.Range("O13").Value = 20.64103213442915
$ws.Range("B14").Value = 17.92890836825293
$ws.Range("C14").Value = 14.45275319775049
$ws.Range("E14").Value = 17.15218754573755
$ws.Range("F14").Value = 35.37645660104224
$ws.Range("G14").Value = 26.67801326970192
$ws.Range("H14").Value = 13.75217105632054
$ws.Range("J14").Value = 7.720424909456193
$ws.Range("L14").Value = 12.77767426947507
$ws.Range("M14").Value = 17.1032130942892
$ws.Range("O14").Value = 20.65780453386639
$ws.Range("B15").Value = 17.86303455565788
$ws.Range("C15").Value = 14.44310066706836
$ws.Range("E15").Value = 17.1577699697289
$ws.Range("F15").Value = 35.38406656922813
$ws.Range("G15").Value = 26.68916051445493
$ws.Range("H15").Value = 13.75840944306252
$ws.Range("J15").Value = 7.718386347224452
$ws.Range("L15").Value = 12.77469722795087
$ws.Range("M15").Value = 17.08308132002171
$ws.Range("O15").Value = 20.66818416688665
$ws.Range("B16").Value = 17.48070136645498
$ws.Range("C16").Value = 14.38794310275688
$ws.Range("E16").Value = 17.1904339096927
$ws.Range("F16").Value = 35.43048911395302
$ws.Range("G16").Value = 26.75643511491248
$ws.Range("H16").Value = 13.79492887895691
$ws.Range("J16").Value = 7.706720004192364
$ws.Range("L16").Value = 12.7580034873221
$ws.Range("M16").Value = 16.9676065909744
$ws.Range("O16").Value = 20.7292886441676
$ws.Range("B17").Value = 17.2419196191434
$ws.Range("C17").Value = 14.35425531236459
$ws.Range("E17").Value = 17.21107348648617
$ws.Range("F17").Value = 35.46147425841721
$ws.Range("G17").Value = 26.80072188787504
$ws.Range("H17").Value = 13.81801908114624
$ws.Range("J17").Value = 7.69957808805857
$ws.Range("L17").Value = 12.74809200754267
$ws.Range("M17").Value = 16.89669564413244
$ws.Range("O17").Value = 20.76821975500494
$ws.Range("B18").Value = 17.10305800184155
$ws.Range("C18").Value = 14.33493764611569
$ws.Range("E18").Value = 17.22316599825544
$ws.Range("F18").Value = 35.48021571067974
$ws.Range("G18").Value = 26.82729769038612
$ws.Range("H18").Value = 13.83155225711579
$ws.Range("J18").Value = 7.69547610882346
$ws.Range("L18").Value = 12.74251331301176
$ws.Range("M18").Value = 16.85589348093995
$ws.Range("O18").Value = 20.79114203239328
$ws.Range("B19").Value = 17.05578414236977
$ws.Range("C19").Value = 14.32840743952713
$ws.Range("E19").Value = 17.2272983202356
$ws.Range("F19").Value = 35.48671904292181
$ws.Range("G19").Value = 26.83648478489569
$ws.Range("H19").Value = 13.83617768559645
$ws.Range("J19").Value = 7.694088288675332
$ws.Range("L19").Value = 12.74064553188851
$ws.Range("M19").Value = 16.84207676374776
$ws.Range("O19").Value = 20.79899407010635
$ws.Range("B20").Value = 17.26749654258399
$ws.Range("C20").Value = 14.35783544700362
$ws.Range("E20").Value = 17.20885348212669
$ws.Range("F20").Value = 35.4580806255388
$ws.Range("G20").Value = 26.79589318869128
$ws.Range("H20").Value = 13.81553497600629
$ws.Range("J20").Value = 7.700337758313747
$ws.Range("L20").Value = 12.74913448421504
$ws.Range("M20").Value = 16.90424612415142
$ws.Range("O20").Value = 20.76402058075967
$ws.Range("B21").Value = 17.96040996716695
$ws.Range("C21").Value = 14.45738449748349
$ws.Range("E21").Value = 17.14952252552896
$ws.Range("F21").Value = 35.37285722085095
$ws.Range("G21").Value = 26.67272808160624
$ws.Range("H21").Value = 13.7491931992754
$ws.Range("J21").Value = 7.721402723092085
$ws.Range("L21").Value = 12.77910828788753
$ws.Range("M21").Value = 17.11286457041911
$ws.Range("O21").Value = 20.65285596595513
$ws.Range("B22").Value = 18.39902176367073
$ws.Range("C22").Value = 14.5228893860584
$ws.Range("E22").Value = 17.11271173134918
$ws.Range("F22").Value = 35.32537811869073
$ws.Range("G22").Value = 26.60216423227839
$ws.Range("H22").Value = 13.70808310797572
$ws.Range("J22").Value = 7.735214692935284
$ws.Range("L22").Value = 12.79976378976082
$ws.Range("M22").Value = 17.24885406087682
$ws.Range("O22").Value = 20.58494727664401
$ws.Range("B23").Value = 18.16626142366944
$ws.Range("C23").Value = 14.48789140095085
$ws.Range("E23").Value = 17.13217870565317
$ws.Range("F23").Value = 35.34996518195012
$ws.Range("G23").Value = 26.63891166703694
$ws.Range("H23").Value = 13.72981856504807
$ws.Range("J23").Value = 7.727839246834513
$ws.Range("L23").Value = 12.78864310360902
$ws.Range("M23").Value = 17.17631653288118
$ws.Range("O23").Value = 20.62075624195183
$ws.Range("B24").Value = 17.25593813739956
$ws.Range("C24").Value = 14.35621671184775
$ws.Range("E24").Value = 17.20985644038869
$ws.Range("F24").Value = 35.45961199766047
$ws.Range("G24").Value = 26.79807277243605
$ws.Range("H24").Value = 13.81665723532079
$ws.Range("J24").Value = 7.699994298833837
$ws.Range("L24").Value = 12.74866280790147
$ws.Range("M24").Value = 16.90083265726703
$ws.Range("O24").Value = 20.76591734479065
$ws.Range("B25").Value = 16.21643387254255
$ws.Range("C25").Value = 14.21647185860149
$ws.Range("E25").Value = 17.3021030467803
$ws.Range("F25").Value = 35.61301944965331
$ws.Range("G25").Value = 27.01186192647678
$ws.Range("H25").Value = 13.91997242807043
$ws.Range("J25").Value = 7.670187262271609
$ws.Range("L25").Value = 12.71021800375828
$ws.Range("M25").Value = 16.60319208306498
$ws.Range("O25").Value = 20.94275191450809
